# Fix values on the "sqlCount" sheet: sqlRecordCount and sqlColCount
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sqlCount")

# Values must stay text (shared strings), matching original cell formatting
$ws.Range("A2").NumberFormat = "@"
$ws.Range("B2").NumberFormat = "@"

$ws.Range("A2").Value = "252"
$ws.Range("B2").Value = "1"
